{"js": "// Lattice multiplication exercises: replace the multiplication problems\n// shown in each table cell with a new set of problems, keeping the same\n// 5-row x 3-column table layout and per-cell formatting (font size 32).\n//\n// Each cell holds a single run made of 5 text segments separated by\n// manual line breaks:\n//   \"AA x BB\"\n//   \"  b1    b2\"\n//   \"  ----\"\n//   \"a1|    |\"\n//   \"a2|    |\"\n// Only the text segments change; the table/row/cell/run structure itself\n// is untouched.\n\nconst newCellLines = [\n  [\"75 x 76\", \"  7    6\", \"  ----\", \"7|    |\", \"5|    |\"],\n  [\"80 x 54\", \"  5    4\", \"  ----\", \"8|    |\", \"0|    |\"],\n  [\"84 x 41\", \"  4    1\", \"  ----\", \"8|    |\", \"4|    |\"],\n  [\"57 x 37\", \"  3    7\", \"  ----\", \"5|    |\", \"7|    |\"],\n  [\"40 x 69\", \"  6    9\", \"  ----\", \"4|    |\", \"0|    |\"],\n  [\"37 x 55\", \"  5    5\", \"  ----\", \"3|    |\", \"7|    |\"],\n  [\"43 x 68\", \"  6    8\", \"  ----\", \"4|    |\", \"3|    |\"],\n  [\"29 x 19\", \"  1    9\", \"  ----\", \"2|    |\", \"9|    |\"],\n  [\"90 x 93\", \"  9    3\", \"  ----\", \"9|    |\", \"0|    |\"],\n  [\"68 x 71\", \"  7    1\", \"  ----\", \"6|    |\", \"8|    |\"],\n  [\"64 x 12\", \"  1    2\", \"  ----\", \"6|    |\", \"4|    |\"],\n  [\"11 x 26\", \"  2    6\", \"  ----\", \"1|    |\", \"1|    |\"],\n  [\"72 x 66\", \"  6    6\", \"  ----\", \"7|    |\", \"2|    |\"],\n  [\"13 x 36\", \"  3    6\", \"  ----\", \"1|    |\", \"3|    |\"],\n  [\"10 x 98\", \"  9    8\", \"  ----\", \"1|    |\", \"0|    |\"],\n];\n\nconst VT = \"\\u000b\"; // vertical tab -> manual line break (<w:br/>) on insertText\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\n\nlet cellIdx = 0;\nfor (let r = 0; r < rows.length; r++) {\n  const row = rows[r];\n  row.cells.load(\"items\");\n  await context.sync();\n\n  const cells = row.cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    if (cellIdx >= newCellLines.length) break;\n    const lines = newCellLines[cellIdx];\n    const text = lines.join(VT);\n\n    const cellRange = cells[c].body.getRange(\"Whole\");\n    cellRange.insertText(text, \"Replace\");\n\n    cellIdx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication exercises: replace the multiplication problems\n# shown in each table cell with a new set of problems, keeping the same\n# 5-row x 3-column table layout and per-cell formatting (font size 32).\n#\n# Each cell holds a single run made of 5 text \"lines\" separated by manual\n# line breaks (vertical-tab char, 0x0B, which Word's COM model renders as\n# <w:br/> between <w:t> runs):\n#   \"AA x BB\"\n#   \"  b1    b2\"\n#   \"  ----\"\n#   \"a1|    |\"\n#   \"a2|    |\"\n# Only the text content changes; the table/row/cell/run structure itself\n# is left alone by assigning straight into Cell(r,c).Range.Text.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$nl = [char]11\n\n$newCellLines = @(\n    @(\"75 x 76\", \"  7    6\", \"  ----\", \"7|    |\", \"5|    |\"),\n    @(\"80 x 54\", \"  5    4\", \"  ----\", \"8|    |\", \"0|    |\"),\n    @(\"84 x 41\", \"  4    1\", \"  ----\", \"8|    |\", \"4|    |\"),\n    @(\"57 x 37\", \"  3    7\", \"  ----\", \"5|    |\", \"7|    |\"),\n    @(\"40 x 69\", \"  6    9\", \"  ----\", \"4|    |\", \"0|    |\"),\n    @(\"37 x 55\", \"  5    5\", \"  ----\", \"3|    |\", \"7|    |\"),\n    @(\"43 x 68\", \"  6    8\", \"  ----\", \"4|    |\", \"3|    |\"),\n    @(\"29 x 19\", \"  1    9\", \"  ----\", \"2|    |\", \"9|    |\"),\n    @(\"90 x 93\", \"  9    3\", \"  ----\", \"9|    |\", \"0|    |\"),\n    @(\"68 x 71\", \"  7    1\", \"  ----\", \"6|    |\", \"8|    |\"),\n    @(\"64 x 12\", \"  1    2\", \"  ----\", \"6|    |\", \"4|    |\"),\n    @(\"11 x 26\", \"  2    6\", \"  ----\", \"1|    |\", \"1|    |\"),\n    @(\"72 x 66\", \"  6    6\", \"  ----\", \"7|    |\", \"2|    |\"),\n    @(\"13 x 36\", \"  3    6\", \"  ----\", \"1|    |\", \"3|    |\"),\n    @(\"10 x 98\", \"  9    8\", \"  ----\", \"1|    |\", \"0|    |\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $newCellLines.Count) { break }\n        $lines = $newCellLines[$idx]\n        $text = ($lines -join $nl)\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $text\n        $idx++\n    }\n}\n"}
